$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 134.8889375
$ws.Range("C3").Value = 50.00226666666666
$ws.Range("C4").Value = 123.84800625
$ws.Range("C5").Value = 17.02913333333333
$ws.Range("C6").Value = 74.68485
$ws.Range("C7").Value = 92.43515000000001
$ws.Range("C9").Value = 0.0309375
$ws.Range("C10").Value = 3.8105
$ws.Range("C12").Value = 0.3396875
$ws.Range("C15").Value = 0.130625
$ws.Range("C16").Value = 0.04575
$ws.Range("C17").Value = 0.1985625
$ws.Range("C18").Value = 0.01375
$ws.Range("C19").Value = 0.3193125
$ws.Range("C20").Value = 0.0909375
$ws.Range("C21").Value = 0.5634666666666667
$ws.Range("C22").Value = 0.239125
$ws.Range("C23").Value = 71.94475
$ws.Range("C24").Value = 0.2743333333333334
$ws.Range("C26").Value = 0.021375
$ws.Range("C27").Value = 1.2431
$ws.Range("C28").Value = 0.1996666666666667
$ws.Range("C30").Value = 0.0665
$ws.Range("C36").Value = 0.2142
$ws.Range("C37").Value = 0.3381666666666667
$ws.Range("C39").Value = 5.749583333333334
$ws.Range("C41").Value = 0.4808333333333333
$ws.Range("C42").Value = 0.004
$ws.Range("C43").Value = 0.0274375
$ws.Range("C45").Value = 0.006562500000000001
$ws.Range("C46").Value = 0.381875
$ws.Range("C47").Value = 0.189
$ws.Range("C48").Value = 15.0102625
$ws.Range("C49").Value = 3.2697875
$ws.Range("C51").Value = 0.1525
$ws.Range("C52").Value = 0.7208333333333333
$ws.Range("C53").Value = 1.027666666666667
$ws.Range("C55").Value = 9.226666666666667
$ws.Range("C56").Value = 0.850375
$ws.Range("C58").Value = 0.080625
$ws.Range("C59").Value = 0.0508125
$ws.Range("C60").Value = 0.117375
$ws.Range("C62").Value = 0.0075
$ws.Range("C64").Value = 0.83
$ws.Range("C65").Value = 0.0121875
$ws.Range("C66").Value = 0.8963333333333333
$ws.Range("C67").Value = 0.1
$ws.Range("C68").Value = 0.0621875
$ws.Range("C69").Value = 0.0140625
$ws.Range("C70").Value = 0.25575
$ws.Range("C72").Value = 0.011875
$ws.Range("C73").Value = 0.0075
$ws.Range("C74").Value = 0.5833333333333334
$ws.Range("C75").Value = 0.023125
$ws.Range("C76").Value = 0.125
$ws.Range("C78").Value = 0.01778125
$ws.Range("C79").Value = 0.23125
$ws.Range("C80").Value = 0.0515
$ws.Range("C82").Value = 0.48675
$ws.Range("C83").Value = 0.019125
$ws.Range("C86").Value = 0.764525
$ws.Range("C87").Value = 4.742833333333333
$ws.Range("C89").Value = 0.1846875
$ws.Range("C90").Value = 2.36319375
$ws.Range("C92").Value = 0.0096875
$ws.Range("C93").Value = 0.2791875
$ws.Range("C94").Value = 0.04566666666666667
$ws.Range("C95").Value = 0.1781666666666667
$ws.Range("C96").Value = 2.91838125
$ws.Range("C97").Value = 6.864
$ws.Range("C98").Value = 0.4834166666666667
$ws.Range("C99").Value = 0.05624999999999999
$ws.Range("C102").Value = 0.00375
$ws.Range("C103").Value = 1.19725
$ws.Range("C104").Value = 0.025
$ws.Range("C106").Value = 0.383625
$ws.Range("C107").Value = 3.478866666666667
$ws.Range("C108").Value = 1.3423125
$ws.Range("C110").Value = 0.362
$ws.Range("C111").Value = 32.73471666666666
$ws.Range("C112").Value = 0.1558125
$ws.Range("C114").Value = 1.7024375
$ws.Range("C115").Value = 0.020625
$ws.Range("C116").Value = 0.104375
$ws.Range("C117").Value = 18.83569375
$ws.Range("C118").Value = 17.56473333333333
$ws.Range("C119").Value = 16.54975625
$ws.Range("C120").Value = 7.7516
$ws.Range("C121").Value = 24.23241666666667
$ws.Range("C122").Value = 2.621125
$ws.Range("C124").Value = 0.0234375
$ws.Range("C125").Value = 1.010666666666667
$ws.Range("C130").Value = 0.02
$ws.Range("C131").Value = 0.00825
$ws.Range("C132").Value = 0.150375
$ws.Range("C136").Value = 0.4951333333333334
$ws.Range("C137").Value = 0.21825
$ws.Range("C138").Value = 0.28625
$ws.Range("C139").Value = 0.03933333333333334
$ws.Range("C142").Value = 0.2211666666666667
$ws.Range("C143").Value = 0.06849999999999999
$ws.Range("C151").Value = 0.2076666666666667
$ws.Range("C152").Value = 0.3178333333333334
$ws.Range("C154").Value = 2.739416666666667
$ws.Range("C156").Value = 0.4808333333333334
$ws.Range("C160").Value = 0.006562500000000001
$ws.Range("C162").Value = 0.1473333333333333
$ws.Range("C163").Value = 0.02625
$ws.Range("C164").Value = 0.1190625
$ws.Range("C166").Value = 0.1525
$ws.Range("C167").Value = 0.7208333333333333
$ws.Range("C168").Value = 0.9716666666666667
$ws.Range("C170").Value = 0.2361333333333333
$ws.Range("C171").Value = 0.0733125
$ws.Range("C173").Value = 0.025625
$ws.Range("C174").Value = 0.0320625
$ws.Range("C179").Value = 0.463125
$ws.Range("C181").Value = 0.2333333333333333
$ws.Range("C184").Value = 0.0140625
$ws.Range("C188").Value = 0.0075
$ws.Range("C189").Value = 0.5833333333333334
$ws.Range("C190").Value = 0.00625
$ws.Range("C191").Value = 0
$ws.Range("C195").Value = 0.0515
$ws.Range("C198").Value = 0.019125
$ws.Range("C201").Value = 0.031375
$ws.Range("C202").Value = 1.246816666666667
$ws.Range("C204").Value = 0.1815625
$ws.Range("C205").Value = 0.2515625
$ws.Range("C208").Value = 0.045625
$ws.Range("C209").Value = 0.011
$ws.Range("C210").Value = 0
$ws.Range("C211").Value = 1.9066875
$ws.Range("C212").Value = 1.4
$ws.Range("C213").Value = 0.1443666666666667
$ws.Range("C218").Value = 0.6851875000000001
$ws.Range("C221").Value = 0.2243125
$ws.Range("C222").Value = 0
$ws.Range("C225").Value = 0
$ws.Range("C226").Value = 12.80455
